$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E. This shifts the old "Is model" column
# (previously E) one position right to F, while the "PDB filename" column
# (D) stays put.
$ws.Columns("E").Insert()

# --- Header row ---
$ws.Range("D1").Value = "PDB or RCSB ID"
$ws.Range("E1").Value = "File Extension"
$ws.Range("F1").Value = " Is model"
$ws.Range("G1").Value = "From RCSB"

# --- New "File Extension" column values (all "pdb") ---
$ws.Range("E2").Value = "pdb"
$ws.Range("E3").Value = "pdb"
$ws.Range("E4").Value = "pdb"
$ws.Range("E5").Value = "pdb"

# --- New "From RCSB" column values (all "n") ---
$ws.Range("G2").Value = "n"
$ws.Range("G3").Value = "n"
$ws.Range("G4").Value = "n"
$ws.Range("G5").Value = "n"

# --- Column widths ---
# C keeps its original width; D (now "PDB or RCSB ID") takes on the same
# width as C; E (the new "File Extension" column) takes on the width that D
# used to have.
$ws.Columns("C").ColumnWidth = 19
$ws.Columns("D").ColumnWidth = 19
$ws.Columns("E").ColumnWidth = 21.5

# --- Selection / active cell ---
$ws.Range("G5").Select()
